$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Marketing & Sales" list in the Internships 2022 sheet (rows 17-30)
# gets a new entry for Tiffany Than at the very top. Her existing entry
# (the last row of that section, row 30) is moved up to row 17 and every
# other row in the section shifts down by one.

# 1) Make room for the moved-up entry by inserting a blank row at 17;
#    this pushes the old rows 17-30 down to 18-31 (Tiffany Than's
#    original row is now row 31).
$ws.Rows("17").Insert()

# 2) Copy Tiffany Than's row (now row 31) into the newly inserted row 17.
$ws.Rows("31").Copy()
$ws.Paste($ws.Rows("17"))

# The paste brings the correct style for columns B-D, but column A keeps
# the "section header" style that was inherited when the blank row was
# inserted (it copied formatting from row 16, the Human Resources
# header). Fix it up to match the plain entry style used by the rest of
# the Marketing & Sales rows (e.g. row 18, formerly row 17).
$ws.Range("A17").Font.Color = $ws.Range("A18").Font.Color
$ws.Range("A17").Font.Bold = $ws.Range("A18").Font.Bold

# 3) Remove the now-duplicated row that the copy left behind.
$ws.Rows("31").Delete()

# Update the saved cursor/selection to match the authored workbook state.
$ws.Range("B14").Select()
